# Pygame.pptx edit: slide 3 ("Реализация") — the bullet list shape that
# ends with "Игрок" gets a new trailing bullet "Сетка" appended as its
# own paragraph.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# Append a new paragraph ("\r" = paragraph break) holding "Сетка" after
# the existing "...Игрок" text, mirroring pressing Enter + typing at the
# end of that placeholder.
$tr.InsertAfter("`rСетка") | Out-Null
